$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 54.86839566666666
$ws.Range("H2").Value = 164.605187
$ws.Range("I2").Value = 0.6170939026906647
$ws.Range("J2").Value = 0.6170939026906647
$ws.Range("M2").Value = 2.598166333333333
$ws.Range("N2").Value = 7.794499
$ws.Range("O2").Value = 0.3466013321552429
$ws.Range("P2").Value = 0.3466013321552429
$ws.Range("Q2").Value = 142.5572183851459
$ws.Range("R2").Value = 1283.014965466313
$ws.Range("S2").Value = 0.2138855687374622
$ws.Range("T2").Value = 0.2138855687374622

$ws.Range("G3").Value = 54.86839566666666
$ws.Range("H3").Value = 164.605187
$ws.Range("I3").Value = 0.6170939026906647
$ws.Range("J3").Value = 0.6170939026906647
$ws.Range("M3").Value = 4.333403333333333
$ws.Range("O3").Value = 0.5780859172985858
$ws.Range("P3").Value = 0.5780859172985858
$ws.Range("Q3").Value = 237.7668886765855
$ws.Range("R3").Value = 2139.90199808927
$ws.Range("S3").Value = 0.3567332947962971
$ws.Range("T3").Value = 0.3567332947962971

$ws.Range("G4").Value = 54.86839566666666
$ws.Range("H4").Value = 164.605187
$ws.Range("I4").Value = 0.6170939026906647
$ws.Range("J4").Value = 0.6170939026906647
$ws.Range("M4").Value = 0.4692043333333333
$ws.Range("N4").Value = 1.407613
$ws.Range("O4").Value = 0.06259293136852516
$ws.Range("P4").Value = 0.06259293136852516
$ws.Range("Q4").Value = 25.74448900984789
$ws.Range("R4").Value = 231.700401088631
$ws.Range("S4").Value = 0.03862571629905211
$ws.Range("T4").Value = 0.03862571629905211

$ws.Range("G5").Value = 54.86839566666666
$ws.Range("H5").Value = 164.605187
$ws.Range("I5").Value = 0.6170939026906647
$ws.Range("J5").Value = 0.6170939026906647
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.09534933333333333
$ws.Range("N5").Value = 0.286048
$ws.Range("O5").Value = 0.01271981917764605
$ws.Range("P5").Value = 0.01271981917764604
$ws.Range("Q5").Value = 5.231664947886221
$ws.Range("R5").Value = 47.08498453097599
$ws.Range("S5").Value = 0.007849322857853159
$ws.Range("T5").Value = 0.007849322857853159

$ws.Range("I6").Value = 0.06720170646055251
$ws.Range("J6").Value = 0.0672017064605525
$ws.Range("M6").Value = 2.598166333333333
$ws.Range("N6").Value = 7.794499
$ws.Range("O6").Value = 0.3466013321552429
$ws.Range("P6").Value = 0.3466013321552429
$ws.Range("Q6").Value = 15.52452277032745
$ws.Range("R6").Value = 139.720704932947
$ws.Range("S6").Value = 0.0232922009823331
$ws.Range("T6").Value = 0.02329220098233309

$ws.Range("I7").Value = 0.06720170646055251
$ws.Range("J7").Value = 0.0672017064605525
$ws.Range("M7").Value = 4.333403333333333
$ws.Range("O7").Value = 0.5780859172985858
$ws.Range("P7").Value = 0.5780859172985858
$ws.Range("Q7").Value = 25.89288370734778
$ws.Range("S7").Value = 0.0388483601232788
$ws.Range("T7").Value = 0.0388483601232788

$ws.Range("I8").Value = 0.06720170646055251
$ws.Range("J8").Value = 0.0672017064605525
$ws.Range("M8").Value = 0.4692043333333333
$ws.Range("N8").Value = 1.407613
$ws.Range("O8").Value = 0.06259293136852516
$ws.Range("P8").Value = 0.06259293136852516
$ws.Range("Q8").Value = 2.803582381665445
$ws.Range("R8").Value = 25.232241434989
$ws.Range("S8").Value = 0.004206351800333137
$ws.Range("T8").Value = 0.004206351800333136

$ws.Range("I9").Value = 0.06720170646055251
$ws.Range("J9").Value = 0.0672017064605525
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.09534933333333333
$ws.Range("N9").Value = 0.286048
$ws.Range("O9").Value = 0.01271981917764605
$ws.Range("P9").Value = 0.01271981917764604
$ws.Range("Q9").Value = 0.5697298427271111
$ws.Range("R9").Value = 5.127568584544
$ws.Range("S9").Value = 0.000854793554607476
$ws.Range("T9").Value = 0.0008547935546074757

$ws.Range("G10").Value = 4.832157666666667
$ws.Range("H10").Value = 14.496473
$ws.Range("I10").Value = 0.05434631351453007
$ws.Range("J10").Value = 0.05434631351453007
$ws.Range("M10").Value = 2.598166333333333
$ws.Range("N10").Value = 7.794499
$ws.Range("O10").Value = 0.3466013321552429
$ws.Range("P10").Value = 0.3466013321552429
$ws.Range("Q10").Value = 12.55474936689189
$ws.Range("R10").Value = 112.992744302027
$ws.Range("S10").Value = 0.0188365046618626
$ws.Range("T10").Value = 0.0188365046618626

$ws.Range("G11").Value = 4.832157666666667
$ws.Range("H11").Value = 14.496473
$ws.Range("I11").Value = 0.05434631351453007
$ws.Range("J11").Value = 0.05434631351453007
$ws.Range("M11").Value = 4.333403333333333
$ws.Range("O11").Value = 0.5780859172985858
$ws.Range("P11").Value = 0.5780859172985858
$ws.Range("Q11").Value = 20.93968813992556
$ws.Range("R11").Value = 188.45719325933
$ws.Range("S11").Value = 0.03141683849984365
$ws.Range("T11").Value = 0.03141683849984365

$ws.Range("G12").Value = 4.832157666666667
$ws.Range("H12").Value = 14.496473
$ws.Range("I12").Value = 0.05434631351453007
$ws.Range("J12").Value = 0.05434631351453007
$ws.Range("M12").Value = 0.4692043333333333
$ws.Range("N12").Value = 1.407613
$ws.Range("O12").Value = 0.06259293136852516
$ws.Range("P12").Value = 0.06259293136852516
$ws.Range("Q12").Value = 2.267269316549889
$ws.Range("R12").Value = 20.405423848949
$ws.Range("S12").Value = 0.003401695071947332
$ws.Range("T12").Value = 0.003401695071947332

$ws.Range("G13").Value = 4.832157666666667
$ws.Range("H13").Value = 14.496473
$ws.Range("I13").Value = 0.05434631351453007
$ws.Range("J13").Value = 0.05434631351453007
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.09534933333333333
$ws.Range("N13").Value = 0.286048
$ws.Range("O13").Value = 0.01271981917764605
$ws.Range("P13").Value = 0.01271981917764604
$ws.Range("Q13").Value = 0.4607430120782223
$ws.Range("R13").Value = 4.146687108704
$ws.Range("S13").Value = 0.0006912752808764841
$ws.Range("T13").Value = 0.0006912752808764839

$ws.Range("G14").Value = 23.238438
$ws.Range("H14").Value = 69.71531400000001
$ws.Range("I14").Value = 0.2613580773342528
$ws.Range("J14").Value = 0.2613580773342528
$ws.Range("M14").Value = 2.598166333333333
$ws.Range("N14").Value = 7.794499
$ws.Range("O14").Value = 0.3466013321552429
$ws.Range("P14").Value = 0.3466013321552429
$ws.Range("Q14").Value = 60.37732725085401
$ws.Range("R14").Value = 543.3959452576861
$ws.Range("S14").Value = 0.090587057773585
$ws.Range("T14").Value = 0.090587057773585

$ws.Range("G15").Value = 23.238438
$ws.Range("H15").Value = 69.71531400000001
$ws.Range("I15").Value = 0.2613580773342528
$ws.Range("J15").Value = 0.2613580773342528
$ws.Range("M15").Value = 4.333403333333333
$ws.Range("O15").Value = 0.5780859172985858
$ws.Range("P15").Value = 0.5780859172985858
$ws.Range("Q15").Value = 100.70152469066
$ws.Range("R15").Value = 906.3137222159401
$ws.Range("S15").Value = 0.1510874238791662
$ws.Range("T15").Value = 0.1510874238791662

$ws.Range("G16").Value = 23.238438
$ws.Range("H16").Value = 69.71531400000001
$ws.Range("I16").Value = 0.2613580773342528
$ws.Range("J16").Value = 0.2613580773342528
$ws.Range("M16").Value = 0.4692043333333333
$ws.Range("N16").Value = 1.407613
$ws.Range("O16").Value = 0.06259293136852516
$ws.Range("P16").Value = 0.06259293136852516
$ws.Range("Q16").Value = 10.903575809498
$ws.Range("R16").Value = 98.132182285482
$ws.Range("S16").Value = 0.01635916819719258
$ws.Range("T16").Value = 0.01635916819719258

$ws.Range("G17").Value = 23.238438
$ws.Range("H17").Value = 69.71531400000001
$ws.Range("I17").Value = 0.2613580773342528
$ws.Range("J17").Value = 0.2613580773342528
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.09534933333333333
$ws.Range("N17").Value = 0.286048
$ws.Range("O17").Value = 0.01271981917764605
$ws.Range("P17").Value = 0.01271981917764604
$ws.Range("Q17").Value = 2.215769571008
$ws.Range("R17").Value = 19.941926139072
$ws.Range("S17").Value = 0.003324427484308927
$ws.Range("T17").Value = 0.003324427484308926
